$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 33471.332
$ws.Range("I31").Value = 33471.332
$ws.Range("K31").Value = 100413.996
$ws.Range("M31").Value = -100183.996
$ws.Range("H51").Value = 7479.1
$ws.Range("I51").Value = 5957.8
$ws.Range("J51").Value = 9000.4
$ws.Range("K51").Value = 5957.8
$ws.Range("L51").Value = 9000.4
$ws.Range("M51").Value = -5473.8
$ws.Range("N51").Value = -9968.4
$ws.Range("H87").Value = 19999.818
$ws.Range("J87").Value = 19999.818
$ws.Range("L87").Value = 19999.818
$ws.Range("N87").Value = -22495.818
$ws.Range("H90").Value = 19999.818
$ws.Range("J90").Value = 19999.818
$ws.Range("L90").Value = 59999.454
$ws.Range("N90").Value = -72479.454
$ws.Range("H96").Value = 228.75
$ws.Range("I96").Value = 152.85715
$ws.Range("K96").Value = 458.57145
$ws.Range("M96").Value = 914.4285500000001
$ws.Range("H113").Value = 6584.4546
$ws.Range("J113").Value = 6581.1113
$ws.Range("L113").Value = 6581.1113
$ws.Range("N113").Value = -13089.1113
$ws.Range("H132").Value = 28569.158
$ws.Range("I132").Value = 1754.9231
$ws.Range("K132").Value = 5264.7693
$ws.Range("M132").Value = -2734.7693
$ws.Range("H138").Value = 2405.015
$ws.Range("I138").Value = 1553.0555
$ws.Range("J138").Value = 2717.9795
$ws.Range("K138").Value = 4659.166499999999
$ws.Range("L138").Value = 8153.9385
$ws.Range("M138").Value = 480.8335000000006
$ws.Range("N138").Value = -18433.9385
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8581.717000000001
$ws.Range("I32").Value = 6984.224
$ws.Range("K32").Value = 6984.224
$ws.Range("M32").Value = -6697.224
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9258
$ws.Range("H97").Value = 1865.3077
$ws.Range("I97").Value = 1727.7778
$ws.Range("K97").Value = 1727.7778
$ws.Range("M97").Value = -1231.7778
$ws.Range("H102").Value = 1929.4
$ws.Range("I102").Value = 1929.4
$ws.Range("K102").Value = 1929.4
$ws.Range("M102").Value = -307.4000000000001
$ws.Range("H132").Value = 6840.4443
$ws.Range("I132").Value = 7310
$ws.Range("J132").Value = 6253.5
$ws.Range("K132").Value = 21930
$ws.Range("L132").Value = 18760.5
$ws.Range("M132").Value = -19400
$ws.Range("N132").Value = -23820.5
$ws.Range("H133").Value = 68189.7
$ws.Range("I133").Value = 63847.4
$ws.Range("K133").Value = 63847.4
$ws.Range("M133").Value = -61317.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 247.3158
$ws.Range("J80").Value = 264.05884
$ws.Range("L80").Value = 264.05884
$ws.Range("N80").Value = -2260.05884
$ws.Range("H83").Value = 247.3158
$ws.Range("J83").Value = 264.05884
$ws.Range("L83").Value = 1320.2942
$ws.Range("N83").Value = -11304.2942
$ws.Range("H86").Value = 6106.1816
$ws.Range("I86").Value = 2260.1667
$ws.Range("K86").Value = 2260.1667
$ws.Range("M86").Value = -1137.1667
$ws.Range("H89").Value = 6106.1816
$ws.Range("I89").Value = 2260.1667
$ws.Range("K89").Value = 11300.8335
$ws.Range("M89").Value = -5684.833500000001
$ws.Range("H107").Value = 591.9167
$ws.Range("I107").Value = 567
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 567
$ws.Range("L107").Value = 666.6667
$ws.Range("M107").Value = 1353
$ws.Range("N107").Value = -4506.6667
$ws.Range("H134").Value = 3821.2122
$ws.Range("I134").Value = 2892.5925
$ws.Range("K134").Value = 8677.7775
$ws.Range("M134").Value = -6142.7775
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4383.727
$ws.Range("J58").Value = 4446.6665
$ws.Range("L58").Value = 4446.6665
$ws.Range("N58").Value = -4852.6665
$ws.Range("H86").Value = 33803.3
$ws.Range("I86").Value = 41255.375
$ws.Range("J86").Value = 3995
$ws.Range("K86").Value = 41255.375
$ws.Range("L86").Value = 3995
$ws.Range("M86").Value = -40132.375
$ws.Range("N86").Value = -6241
$ws.Range("H89").Value = 33803.3
$ws.Range("I89").Value = 41255.375
$ws.Range("J89").Value = 3995
$ws.Range("K89").Value = 206276.875
$ws.Range("L89").Value = 19975
$ws.Range("M89").Value = -200660.875
$ws.Range("N89").Value = -31207
$ws.Range("H105").Value = 13163271
$ws.Range("I105").Value = 1603.4546
$ws.Range("K105").Value = 1603.4546
$ws.Range("M105").Value = 143.5454
$ws.Range("H136").Value = 4383.727
$ws.Range("J136").Value = 4446.6665
$ws.Range("L136").Value = 13339.9995
$ws.Range("N136").Value = -18439.9995
$ws.Range("H141").Value = 90800
$ws.Range("J141").Value = 90800
$ws.Range("L141").Value = 90800
$ws.Range("N141").Value = -101160
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 333.19232
$ws.Range("I60").Value = 375.375
$ws.Range("K60").Value = 1126.125
$ws.Range("M60").Value = -875.125
$ws.Range("H75").Value = 3103.75
$ws.Range("I75").Value = 1097.5
$ws.Range("J75").Value = 3505
$ws.Range("K75").Value = 3292.5
$ws.Range("L75").Value = 10515
$ws.Range("M75").Value = -2294.5
$ws.Range("N75").Value = -12511
$ws.Range("H78").Value = 3103.75
$ws.Range("I78").Value = 1097.5
$ws.Range("J78").Value = 3505
$ws.Range("K78").Value = 9877.5
$ws.Range("L78").Value = 31545
$ws.Range("M78").Value = -4885.5
$ws.Range("N78").Value = -41529
$ws.Range("H81").Value = 2752.5
$ws.Range("J81").Value = 2752.5
$ws.Range("L81").Value = 8257.5
$ws.Range("N81").Value = -10503.5
$ws.Range("H84").Value = 2752.5
$ws.Range("J84").Value = 2752.5
$ws.Range("L84").Value = 24772.5
$ws.Range("N84").Value = -36004.5
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 3000
$ws.Range("M98").Value = -1502
$ws.Range("H131").Value = 6451.9375
$ws.Range("J131").Value = 8748.333000000001
$ws.Range("L131").Value = 26244.999
$ws.Range("N131").Value = -36324.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 7501.5
$ws.Range("I10").Value = 7501.5
$ws.Range("K10").Value = 7501.5
$ws.Range("M10").Value = -7332.5
$ws.Range("H11").Value = 1435700
$ws.Range("I11").Value = 1889000
$ws.Range("J11").Value = 378000
$ws.Range("K11").Value = 1889000
$ws.Range("L11").Value = 378000
$ws.Range("M11").Value = -1888861
$ws.Range("N11").Value = -378278
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1771.4546
$ws.Range("J46").Value = 2557
$ws.Range("L46").Value = 2557
$ws.Range("N46").Value = -2933
$ws.Range("H82").Value = 2342.8125
$ws.Range("I82").Value = 2007.0834
$ws.Range("J82").Value = 3350
$ws.Range("K82").Value = 2007.0834
$ws.Range("L82").Value = 3350
$ws.Range("M82").Value = -1646.0834
$ws.Range("N82").Value = -4072
$ws.Range("H85").Value = 2342.8125
$ws.Range("I85").Value = 2007.0834
$ws.Range("J85").Value = 3350
$ws.Range("K85").Value = 2007.0834
$ws.Range("L85").Value = 3350
$ws.Range("M85").Value = -759.0834
$ws.Range("N85").Value = -5846
$ws.Range("H93").Value = 1113454.9
$ws.Range("I93").Value = 2682.5
$ws.Range("J93").Value = 3334999.8
$ws.Range("K93").Value = 2682.5
$ws.Range("L93").Value = 3334999.8
$ws.Range("M93").Value = -1434.5
$ws.Range("N93").Value = -3337495.8
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H14").Value = 1521.7142
$ws.Range("I14").Value = 720
$ws.Range("K14").Value = 720
$ws.Range("M14").Value = -552
$ws.Range("H132").Value = 2890.8928
$ws.Range("I132").Value = 1774.9474
$ws.Range("J132").Value = 5246.778
$ws.Range("K132").Value = 5324.8422
$ws.Range("L132").Value = 15740.334
$ws.Range("M132").Value = -2794.8422
$ws.Range("N132").Value = -20800.334
$ws.Range("H136").Value = 3077.1333
$ws.Range("I136").Value = 1841.2858
$ws.Range("J136").Value = 5960.778
$ws.Range("K136").Value = 5523.857400000001
$ws.Range("L136").Value = 17882.334
$ws.Range("M136").Value = -2973.857400000001
$ws.Range("N136").Value = -22982.334
